# Update the build timestamp embedded in the "Version" strings from
# "January 30 2026 16.19.47 EST" to "February 02 2026 12.49.33 EST"
# across the "About" sheet and the "Boundaries and methane sources" sheet.

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Boundaries and methane sources")

# A2: "Version: mines - January 30 (built on January 30 2026 16.19.47 EST)"
$cellA2 = $wsAbout.Range("A2")
$cellA2.Value2 = $cellA2.Value2.Replace($oldStamp, $newStamp)

# A6: Recommended citation text referencing the same build version string.
$cellA6 = $wsAbout.Range("A6")
$cellA6.Value2 = $cellA6.Value2.Replace($oldStamp, $newStamp)

# S2:S7 on the data sheet carry the same "mines - January 30 (built on ...)" value.
for ($row = 2; $row -le 7; $row++) {
    $cell = $wsData.Cells.Item($row, 19)  # column S
    $cell.Value2 = $cell.Value2.Replace($oldStamp, $newStamp)
}
